# Mohammad Nabi.xlsx - scrapped IPL batting stats update
# 1. Rename the sheet from "Sheet1" to "Mohammad Nabi"
# 2. Insert a new leading "matchNo" column
# 3. Add two more match rows of batting data (rows 3 and 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data -----------------------------------------------------------------
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

$rows = @(
    @("28th","Sunrisers Hyderabad","Mohammad Nabi","c Anuj Rawat b Mustafizur Rahman","17","5","1","2","340.00","Rajasthan Royals","Delhi","May 02","Royals won by 55 runs"),
    @("55th","Sunrisers Hyderabad","Mohammad Nabi","c Pollard b Chawla","3","4","0","0","75.00","Mumbai Indians","Abu Dhabi","October 08","Mumbai won by 42 runs"),
    @("3rd","Sunrisers Hyderabad","Mohammad Nabi","c Morgan b Prasidh Krishna","14","11","2","0","127.27","Kolkata Knight Riders","Chennai","April 11","KKR won by 10 runs")
)

$colCount = $headers.Length
$rowCount = $rows.Length

# Make sure every touched cell is stored as TEXT so that number-looking
# values ("17", "340.00", ...) round-trip as strings, matching the
# numberStoredAsText semantics of the source sheet.
$lastColLetter = [char](64 + $colCount)
$dataRange = "A1:" + $lastColLetter + (1 + $rowCount)
$ws.Range($dataRange).NumberFormat = "@"

# --- header row -------------------------------------------------------------
for ($col = 1; $col -le $colCount; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# --- data rows ----------------------------------------------------------
for ($r = 0; $r -lt $rowCount; $r++) {
    $rowData = $rows[$r]
    for ($col = 1; $col -le $colCount; $col++) {
        $ws.Cells.Item($r + 2, $col).Value = $rowData[$col - 1]
    }
}

# --- sheet title ----------------------------------------------------------
$ws.Name = "Mohammad Nabi"
